$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column A (requirement IDs) for all new rows 2..11
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# Fill in column B (requirement descriptions). The order below reproduces
# the shared-string table insertion order from the original authoring
# session (strings were typed bottom-up, row 8 first).
$ws.Range("B8").Value = "Integration with OS Calendar"
$ws.Range("B7").Value = "Usage of real personal data"
$ws.Range("B6").Value = "Integration with a geolocalization provider"
$ws.Range("B4").Value = "Implementation of a built-in chat system"
$ws.Range("B3").Value = "Implementation of a complex Search Functionality"
$ws.Range("B5").Value = "Implementation of a complete rating system"
$ws.Range("B2").Value = "Implementation of a workflow for ""User Profile"" functionalities"

# Column B needs to widen to fit the longest description (best-fit style),
# and the selection moves to B10 to match where the author finished editing.
$ws.Columns("B").ColumnWidth = 57.6
$ws.Range("B10").Select()
